$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Duplicate slide 1 -> creates slide 2 right after it
$range = $s1.Duplicate()
$s2 = $p.Slides.Item(2)

# Remove all shapes from the duplicated slide except the "Grupo 7" group
# (shape 1) and the bottom banner picture (originally shape 6, "0 Imagen").
# Delete from the end backwards so indices remain stable.
$s2.Shapes.Item(9).Delete()
$s2.Shapes.Item(8).Delete()
$s2.Shapes.Item(7).Delete()
$s2.Shapes.Item(5).Delete()
$s2.Shapes.Item(4).Delete()
$s2.Shapes.Item(3).Delete()
$s2.Shapes.Item(2).Delete()
